$d = $word.ActiveDocument

# Grab the last paragraph in the document (ends with the useLocation text)
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range

# Insert a blank paragraph after the last paragraph
$r.InsertParagraphAfter()

# Re-fetch the last paragraph (now the new blank one) and insert another
# paragraph after it for the new "Font Awesome" text.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertParagraphAfter()

# The final paragraph now holds the new text.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Font Awesome : provide good and free icons."
